$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 15: add Fails flag ---
$ws.Range("M15").Value = 0

# --- Row 16: add Fails flag ---
$ws.Range("M16").Value = 0

# --- Row 18: add Fails flag ---
$ws.Range("M18").Value = 0

# --- Row 20: complete the row (test 5 / c) with results ---
# Write order matters for shared-string allocation order, so columns are
# written in the exact sequence needed to reproduce the target layout.
$ws.Range("G20").Value = "02:25:51h"
$ws.Range("H20").NumberFormat = "h:mm"
$ws.Range("H20").Value = "01:56m"
$ws.Range("I20").Value = "00:02:49h"
$ws.Range("K20").Value = "02:22:54h"
$ws.Range("J20").Value = "02:23:47h"
$ws.Range("L20").Value = "2,71GB"
$ws.Range("M20").Value = 1

# --- Row 22: new test 6 / a ---
$ws.Range("A22").Value = 6
$ws.Range("B22").Value = "a"
$ws.Range("C22").Value = 337
$ws.Range("D22").Value = 500
$ws.Range("E22").Value = 500
$ws.Range("F22").Value = 144
$ws.Range("H22").NumberFormat = "h:mm"
$ws.Range("H22").Value = "02:14m"
$ws.Range("I22").Value = "00:03:58h"
$ws.Range("K22").Value = "05:50:16h"
$ws.Range("J22").Value = "05:52:03h"
$ws.Range("L22").Value = "2,72GB"
$ws.Range("G22").Value = "05:54:23h"
$ws.Range("M22").Value = 1

# --- Row 23: test 6 / b ---
$ws.Range("B23").Value = "b"
$ws.Range("C23").Value = 337
$ws.Range("D23").Value = 500
$ws.Range("E23").Value = 500
$ws.Range("F23").Value = 144
$ws.Range("G23").Value = "05:56:52h"
$ws.Range("I23").Value = "00:04:40h"
$ws.Range("H23").NumberFormat = "h:mm"
$ws.Range("H23").Value = "01:52m"
$ws.Range("J23").NumberFormat = "h:mm:ss"
$ws.Range("J23").Value = "05:54:54h"
$ws.Range("K23").Value = "05:52:03h"
$ws.Range("L23").Value = "2,72GB"
$ws.Range("M23").Value = 0

# --- Row 24: test 6 / c ---
$ws.Range("B24").Value = "c"
$ws.Range("C24").Value = 337
$ws.Range("D24").Value = 500
$ws.Range("E24").Value = 500
$ws.Range("F24").Value = 144
$ws.Range("M24").Value = 1

# --- Row 26: new test 7 / a (still running, no results yet) ---
$ws.Range("A26").Value = 7
$ws.Range("B26").Value = "a"
$ws.Range("C26").Value = 337
$ws.Range("D26").Value = 2000
$ws.Range("E26").Value = 50
$ws.Range("F26").Value = 144

# --- Row 27: test 7 / b ---
$ws.Range("B27").Value = "b"
$ws.Range("C27").Value = 337
$ws.Range("D27").Value = 2000
$ws.Range("E27").Value = 50
$ws.Range("F27").Value = 144

# --- Row 28: test 7 / c ---
$ws.Range("B28").Value = "c"
$ws.Range("C28").Value = 337
$ws.Range("D28").Value = 2000
$ws.Range("E28").Value = 50
$ws.Range("F28").Value = 144

# --- Row 31: stray formatted (empty) cells ---
$ws.Range("E31").NumberFormat = "h:mm:ss"
$ws.Range("F31").NumberFormat = "h:mm:ss"
$ws.Range("G31").NumberFormat = "h:mm:ss"

# --- Column F width ---
$ws.Columns.Item(6).ColumnWidth = 13

# --- Page setup ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Selection ---
[void]$ws.Range("E29").Select()
